$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.359.02"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.02"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  +1.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.59"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5348"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2658"
$ws.Range("E8").Value = "  +2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06395"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.92"
$ws.Range("E10").Value = "  +2.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07855"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.565"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.678.25"
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.896.02"
$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5544"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8172"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.13"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.381.59"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.674"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.99"
$ws.Range("E21").Value = "  +2.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.27"
$ws.Range("E22").Value = "  +2.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.045"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.44"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.240"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.22"
$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.580"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.288"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9706"
$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.832"
$ws.Range("E36").Value = "  +1.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.427"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5823"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.078.55"
$ws.Range("E40").Value = "  +4.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8636"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.864"
$ws.Range("E42").Value = "  +2.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.32"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.805.73"
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.01"
$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.014"
$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -6.48%  "

$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.071"
$ws.Range("E50").Value = "  +3.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05168"
$ws.Range("E51").Value = "  +0.54%  "
